# Regenerate save_data to use K (strikeouts) instead of Strike# for column G.
# Only column G (header "K") values change, for data rows 3 through 46.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$newK = @{
    3  = 1
    4  = 0
    5  = 1
    6  = 1
    7  = 0
    8  = 1
    9  = 2
    10 = 1
    11 = 1
    12 = 2
    13 = 3
    14 = 4
    15 = 0
    16 = 3
    17 = 0
    18 = 1
    19 = 5
    20 = 4
    21 = 2
    22 = 5
    23 = 2
    24 = 5
    25 = 4
    26 = 2
    27 = 3
    28 = 0
    29 = 3
    30 = 3
    31 = 6
    32 = 5
    33 = 5
    34 = 5
    35 = 4
    36 = 5
    37 = 7
    38 = 10
    39 = 2
    40 = 5
    41 = 4
    42 = 2
    43 = 4
    44 = 7
    45 = 2
    46 = 1
}

foreach ($row in $newK.Keys) {
    $ws.Range("G$row").Value = $newK[$row]
}
